$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.678113579750061
$ws.Range("B1").Value = 2.837880611419678
$ws.Range("C1").Value = 1.738970279693604
$ws.Range("D1").Value = 1.191162824630737
$ws.Range("E1").Value = 0.9817848205566406
